{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same changes as the target OOXML diff:\n//   1. Update the \"Currently working as Lead Developer & Scrum Master (...)\"\n//      bullet to the new tech list.\n//   2. Rename \"Freelance Full Stack Java Developer\" -> \"Freelance Java Developer\"\n//      in the career-overview table.\n//   3. Update the \"Java, WSO2, Docker, REST, OAuth, Ansible, Linux\" techniques\n//      cell to add \"Java 8, Spring Boot\".\n//   4. Remove the three \"famous quotes\" paragraphs near the end of the CV.\n\nconst body = context.document.body;\n\n// --- 1. Summary bullet: tech list update -------------------------------\nconst leadResults = body.search(\n  \"Lead Developer & Scrum Master (Java, WSO2, Docker, REST,\",\n  { matchCase: true }\n);\nleadResults.load(\"items\");\nawait context.sync();\n\nif (leadResults.items.length > 0) {\n  leadResults.items[0].insertText(\n    \"Lead Developer & Scrum Master (Java 8, Spring boot, Docker, REST API\\u2019s,\",\n    \"Replace\"\n  );\n}\n\nconst oauthResults = body.search(\n  \"OAuth, Ansible, Linux, Open Source, JIRA, Git, Jenkins, Infrastructure as code)\",\n  { matchCase: true }\n);\noauthResults.load(\"items\");\nawait context.sync();\n\nif (oauthResults.items.length > 0) {\n  oauthResults.items[0].insertText(\n    \"OAuth, Docker, Ansible, Linux, Open Source, JIRA, Git, Jenkins, Infrastructure as code)\",\n    \"Replace\"\n  );\n}\n\n// --- 2. Career overview table: drop \"Full Stack\" qualifier -------------\nconst roleResults = body.search(\"Freelance Full Stack Java Developer\", {\n  matchCase: true\n});\nroleResults.load(\"items\");\nawait context.sync();\n\nif (roleResults.items.length > 0) {\n  roleResults.items[0].insertText(\"Freelance Java Developer\", \"Replace\");\n}\n\n// --- 3. Techniques cell: add Java 8, Spring Boot ------------------------\nconst techResults = body.search(\n  \"Java, WSO2, Docker, REST, OAuth, Ansible, Linux\",\n  { matchCase: true }\n);\ntechResults.load(\"items\");\nawait context.sync();\n\nif (techResults.items.length > 0) {\n  techResults.items[0].insertText(\n    \"Java 8, Spring Boot, WSO2, Docker, REST, OAuth, Ansible, Linux\",\n    \"Replace\"\n  );\n}\n\n// --- 4. Remove the three quote paragraphs -------------------------------\nconst quoteStarts = [\n  \"\\u201cEverything should be as simple as it can be, but not simpler.\",\n  \"\\u201cPerfection is reached, not when there is no longer anything to add, but\",\n  \"\\u201cThere are two ways of constructing a software design: One way is to\"\n];\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (quoteStarts.some((s) => t.indexOf(s) >= 0)) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Applies the same changes as the target OOXML diff:\n#   1. Update the \"Currently working as Lead Developer & Scrum Master (...)\"\n#      bullet to the new tech list.\n#   2. Rename \"Freelance Full Stack Java Developer\" -> \"Freelance Java Developer\"\n#      in the career-overview table.\n#   3. Update the \"Java, WSO2, Docker, REST, OAuth, Ansible, Linux\" techniques\n#      cell to add \"Java 8, Spring Boot\".\n#   4. Remove the three \"famous quotes\" paragraphs near the end of the CV.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# --- 1. Summary bullet: tech list update --------------------------------\n$find1 = $d.Content.Find\n$find1.Execute(\n    \"Lead Developer & Scrum Master (Java, WSO2, Docker, REST,\",\n    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Lead Developer & Scrum Master (Java 8, Spring boot, Docker, REST API\u2019s,\",\n    $wdReplaceAll\n)\n\n$find2 = $d.Content.Find\n$find2.Execute(\n    \"OAuth, Ansible, Linux, Open Source, JIRA, Git, Jenkins, Infrastructure as code)\",\n    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"OAuth, Docker, Ansible, Linux, Open Source, JIRA, Git, Jenkins, Infrastructure as code)\",\n    $wdReplaceAll\n)\n\n# --- 2. Career overview table: drop \"Full Stack\" qualifier --------------\n$find3 = $d.Content.Find\n$find3.Execute(\n    \"Freelance Full Stack Java Developer\",\n    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Freelance Java Developer\",\n    $wdReplaceAll\n)\n\n# --- 3. Techniques cell: add Java 8, Spring Boot -------------------------\n$find4 = $d.Content.Find\n$find4.Execute(\n    \"Java, WSO2, Docker, REST, OAuth, Ansible, Linux\",\n    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,\n    \"Java 8, Spring Boot, WSO2, Docker, REST, OAuth, Ansible, Linux\",\n    $wdReplaceAll\n)\n\n# --- 4. Remove the three quote paragraphs --------------------------------\n$quoteMarkers = @(\n    \"*Everything should be as simple as it can be, but not simpler.*\",\n    \"*Perfection is reached, not when there is no longer anything to add, but*\",\n    \"*There are two ways of constructing a software design: One way is to*\"\n)\n\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    foreach ($marker in $quoteMarkers) {\n        if ($t -like $marker) {\n            $targets += $p\n            break\n        }\n    }\n}\n\nfor ($i = $targets.Count - 1; $i -ge 0; $i--) {\n    $targets[$i].Range.Delete()\n}\n"}
